$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value2 = 1171783.9
$ws.Range("I17").Value2 = 0
$ws.Range("J17").Value2 = 1171783.9
$ws.Range("K17").Value2 = 0
$ws.Range("L17").Value2 = 3515351.7
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value2 = -3515687.7
# Row 41
$ws.Range("H41").Value2 = 66.666664
$ws.Range("J41").Value2 = 66.666664
$ws.Range("L41").Value2 = 66.666664
$ws.Range("N41").Value2 = -946.666664
# Row 74
$ws.Range("H74").Value2 = 3676.0688
$ws.Range("I74").Value2 = 2605.4614
$ws.Range("K74").Value2 = 2605.4614
$ws.Range("M74").Value2 = -1669.4614
# Row 77
$ws.Range("H77").Value2 = 3676.0688
$ws.Range("I77").Value2 = 2605.4614
$ws.Range("K77").Value2 = 13027.307
$ws.Range("M77").Value2 = -8347.307000000001
# Row 96
$ws.Range("H96").Value2 = 1247.4546
$ws.Range("I96").Value2 = 1242.2858
$ws.Range("J96").Value2 = 1256.5
$ws.Range("K96").Value2 = 3726.8574
$ws.Range("L96").Value2 = 3769.5
$ws.Range("M96").Value2 = -2353.8574
$ws.Range("N96").Value2 = -6515.5
# Row 113
$ws.Range("H113").Value2 = 6734.9287
$ws.Range("I113").Value2 = 4424.5
$ws.Range("K113").Value2 = 4424.5
$ws.Range("M113").Value2 = -1170.5
# Row 116
$ws.Range("H116").Value2 = 2745.125
$ws.Range("I116").Value2 = 2150
$ws.Range("J116").Value2 = 2943.5
$ws.Range("K116").Value2 = 2150
$ws.Range("L116").Value2 = 2943.5
$ws.Range("M116").Value2 = 1292
$ws.Range("N116").Value2 = -9827.5
# Row 131
$ws.Range("H131").Value2 = 4893.353
$ws.Range("I131").Value2 = 3499.75
$ws.Range("K131").Value2 = 10499.25
$ws.Range("M131").Value2 = -5459.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 8915.178
$ws.Range("I32").Value2 = 3833.8235
$ws.Range("K32").Value2 = 3833.8235
$ws.Range("M32").Value2 = -3546.8235
# Row 33
$ws.Range("H33").Value2 = 2000
$ws.Range("I33").Value2 = 2000
$ws.Range("K33").Value2 = 2000
$ws.Range("M33").Value2 = -1671
# Row 61
$ws.Range("H61").Value2 = 100005800
$ws.Range("I61").Value2 = 250001000
$ws.Range("J61").Value2 = 9000
$ws.Range("K61").Value2 = 250001000
$ws.Range("L61").Value2 = 9000
$ws.Range("M61").Value2 = -250000788
$ws.Range("N61").Value2 = -9424
# Row 96
$ws.Range("H96").Value2 = 0
$ws.Range("J96").Value2 = 0
$ws.Range("L96").Value2 = 0
$ws.Range("N96").ClearContents()
# Row 97
$ws.Range("H97").Value2 = 1280.5555
$ws.Range("I97").Value2 = 1220.4
$ws.Range("K97").Value2 = 1220.4
$ws.Range("M97").Value2 = -724.4000000000001
# Row 136
$ws.Range("H136").Value2 = 100005800
$ws.Range("I136").Value2 = 250001000
$ws.Range("J136").Value2 = 9000
$ws.Range("K136").Value2 = 750003000
$ws.Range("L136").Value2 = 27000
$ws.Range("M136").Value2 = -750000450
$ws.Range("N136").Value2 = -32100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value2 = 1477.5333
$ws.Range("I94").Value2 = 845.2857
$ws.Range("J94").Value2 = 2030.75
$ws.Range("K94").Value2 = 845.2857
$ws.Range("L94").Value2 = 2030.75
$ws.Range("M94").Value2 = -394.2857
$ws.Range("N94").Value2 = -2932.75
# Row 105
$ws.Range("H105").Value2 = 25577.25
$ws.Range("I105").Value2 = 25577.25
$ws.Range("K105").Value2 = 25577.25
$ws.Range("M105").Value2 = -23830.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value2 = 8877.25
$ws.Range("I23").Value2 = 8877.25
$ws.Range("J23").Value2 = 0
$ws.Range("K23").Value2 = 8877.25
$ws.Range("L23").Value2 = 0
$ws.Range("M23").Value2 = -8637.25
$ws.Range("N23").ClearContents()
# Row 27
$ws.Range("H27").Value2 = 8877.25
$ws.Range("I27").Value2 = 8877.25
$ws.Range("J27").Value2 = 0
$ws.Range("K27").Value2 = 8877.25
$ws.Range("L27").Value2 = 0
$ws.Range("M27").Value2 = -8685.25
$ws.Range("N27").ClearContents()
# Row 31
$ws.Range("H31").Value2 = 4657.3486
$ws.Range("I31").Value2 = 2511.147
$ws.Range("J31").Value2 = 12765.223
$ws.Range("K31").Value2 = 2511.147
$ws.Range("L31").Value2 = 12765.223
$ws.Range("M31").Value2 = -2216.147
$ws.Range("N31").Value2 = -13355.223
# Row 34
$ws.Range("H34").Value2 = 4657.3486
$ws.Range("I34").Value2 = 2511.147
$ws.Range("J34").Value2 = 12765.223
$ws.Range("K34").Value2 = 2511.147
$ws.Range("L34").Value2 = 12765.223
$ws.Range("M34").Value2 = -2309.147
$ws.Range("N34").Value2 = -13169.223
# Row 86
$ws.Range("H86").Value2 = 4604.273
$ws.Range("I86").Value2 = 4125
$ws.Range("K86").Value2 = 4125
$ws.Range("M86").Value2 = -3002
# Row 89
$ws.Range("H89").Value2 = 4604.273
$ws.Range("I89").Value2 = 4125
$ws.Range("K89").Value2 = 20625
$ws.Range("M89").Value2 = -15009
# Row 122
$ws.Range("H122").Value2 = 1540.7273
$ws.Range("I122").Value2 = 1366.5555
$ws.Range("K122").Value2 = 4099.666499999999
$ws.Range("M122").Value2 = -1649.666499999999
# Row 132
$ws.Range("H132").Value2 = 5595.9033
$ws.Range("I132").Value2 = 5315.7666
$ws.Range("J132").Value2 = 14000
$ws.Range("K132").Value2 = 15947.2998
$ws.Range("L132").Value2 = 42000
$ws.Range("M132").Value2 = -13417.2998
$ws.Range("N132").Value2 = -47060
# Row 134
$ws.Range("H134").Value2 = 1269.0834
$ws.Range("I134").Value2 = 1111.7273
$ws.Range("J134").Value2 = 3000
$ws.Range("K134").Value2 = 3335.1819
$ws.Range("L134").Value2 = 9000
$ws.Range("M134").Value2 = -800.1819
$ws.Range("N134").Value2 = -14070

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 104
$ws.Range("H104").Value2 = 4000
$ws.Range("J104").Value2 = 4000
$ws.Range("L104").Value2 = 12000
$ws.Range("N104").Value2 = -17242

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value2 = 3125379
$ws.Range("I2").Value2 = 7143010
$ws.Range("J2").Value2 = 555.1111
$ws.Range("K2").Value2 = 7143010
$ws.Range("L2").Value2 = 555.1111
$ws.Range("M2").Value2 = -7142897
$ws.Range("N2").Value2 = -781.1111
# Row 51
$ws.Range("H51").Value2 = 79449.75
$ws.Range("I51").Value2 = 79400
$ws.Range("J51").Value2 = 79499.5
$ws.Range("K51").Value2 = 79400
$ws.Range("L51").Value2 = 79499.5
$ws.Range("M51").Value2 = -78891
$ws.Range("N51").Value2 = -80517.5
# Row 122
$ws.Range("H122").Value2 = 3972.5454
$ws.Range("I122").Value2 = 2837.25
$ws.Range("J122").Value2 = 7000
$ws.Range("K122").Value2 = 8511.75
$ws.Range("L122").Value2 = 21000
$ws.Range("M122").Value2 = -6061.75
$ws.Range("N122").Value2 = -25900
# Row 126
$ws.Range("H126").Value2 = 6069.0625
$ws.Range("I126").Value2 = 6551.75
$ws.Range("J126").Value2 = 5586.375
$ws.Range("K126").Value2 = 19655.25
$ws.Range("L126").Value2 = 16759.125
$ws.Range("M126").Value2 = -17185.25
$ws.Range("N126").Value2 = -21699.125
# Row 132
$ws.Range("H132").Value2 = 4146.3
$ws.Range("I132").Value2 = 5245
$ws.Range("K132").Value2 = 15735
$ws.Range("M132").Value2 = -13205

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value2 = 4279.75
$ws.Range("J7").Value2 = 4497.6665
$ws.Range("L7").Value2 = 4497.6665
$ws.Range("N7").Value2 = -4721.6665
# Row 40
$ws.Range("H40").Value2 = 3779.75
$ws.Range("I40").Value2 = 3755.7
$ws.Range("K40").Value2 = 3755.7
$ws.Range("M40").Value2 = -3619.7
# Row 61
$ws.Range("H61").Value2 = 4138.2173
$ws.Range("I61").Value2 = 3640.3333
$ws.Range("J61").Value2 = 4681.364
$ws.Range("K61").Value2 = 3640.3333
$ws.Range("L61").Value2 = 4681.364
$ws.Range("M61").Value2 = -3438.3333
$ws.Range("N61").Value2 = -5085.364
# Row 93
$ws.Range("H93").Value2 = 1843.7241
$ws.Range("I93").Value2 = 856.7143
$ws.Range("J93").Value2 = 2764.9333
$ws.Range("K93").Value2 = 856.7143
$ws.Range("L93").Value2 = 2764.9333
$ws.Range("M93").Value2 = 391.2857
$ws.Range("N93").Value2 = -5260.933300000001
# Row 113
$ws.Range("H113").Value2 = 4138.2173
$ws.Range("I113").Value2 = 3640.3333
$ws.Range("J113").Value2 = 4681.364
$ws.Range("K113").Value2 = 3640.3333
$ws.Range("L113").Value2 = 4681.364
$ws.Range("M113").Value2 = -1470.3333
$ws.Range("N113").Value2 = -9021.364
# Row 126
$ws.Range("H126").Value2 = 4279.75
$ws.Range("J126").Value2 = 4497.6665
$ws.Range("L126").Value2 = 13492.9995
$ws.Range("N126").Value2 = -18432.9995
# Row 132
$ws.Range("H132").Value2 = 1241.7142
$ws.Range("I132").Value2 = 1115.3334
$ws.Range("K132").Value2 = 3346.0002
$ws.Range("M132").Value2 = -816.0001999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value2 = 24001
$ws.Range("I4").Value2 = 47002
$ws.Range("K4").Value2 = 47002
$ws.Range("M4").Value2 = -46889
# Row 22
$ws.Range("H22").Value2 = 7003
$ws.Range("I22").Value2 = 2006.5
$ws.Range("K22").Value2 = 2006.5
$ws.Range("M22").Value2 = -1713.5
# Row 28
$ws.Range("H28").Value2 = 31999.5
# Row 122
$ws.Range("H122").Value2 = 127481.625
$ws.Range("I122").Value2 = 168659
$ws.Range("K122").Value2 = 505977
$ws.Range("M122").Value2 = -503527
# Row 136
$ws.Range("H136").Value2 = 1808.4828
$ws.Range("I136").Value2 = 1722.4615
$ws.Range("K136").Value2 = 5167.3845
$ws.Range("M136").Value2 = -2617.3845
